$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("3rd party"): fill in missing Author for inherits/2.0.3 (row 16, col G) ---
$ws1 = $wb.Worksheets.Item("3rd party")
$ws1.Range("G16").Value = "isaacs"

# --- Sheet 2 ("No IP 3rd party"): add dev-dependency rows ---
$ws2 = $wb.Worksheets.Item("No IP 3rd party")

$devDeps = @(
    @(1, "should/13.2.3", "test framework agnostic BDD-style assertions", "MIT", "https://github.com/shouldjs/should.js", "TJ Holowaychuk"),
    @(2, "should-equal/2.0.0", "Deep comparison of 2 instances for should.js", "MIT", "https://github.com/shouldjs/equal", "Denis Bardadym"),
    @(3, "should-format/3.0.3", "Formatting of objects for should.js", "MIT", "https://github.com/shouldjs/format#readme", "Denis Bardadym"),
    @(4, "should-type/1.4.0", "Simple module to get instance type. Like a bit more advanced version of typeof", "MIT", "https://github.com/shouldjs/type", "Denis Bardadym"),
    @(5, "should-type-adaptors/1.1.0", "Small utility functions to use the same traversing etc code on different types", "MIT", "https://github.com/shouldjs/type-adaptors#readme", "Denis Bardadym"),
    @(6, "should-util/1.0.0", "Utility functions", "MIT", "https://github.com/shouldjs/util#readme", "Denis Bardadym")
)

$rowIndex = 2
foreach ($dep in $devDeps) {
    $ws2.Cells.Item($rowIndex, 1).Value = $dep[0]
    $ws2.Cells.Item($rowIndex, 2).Value = $dep[1]
    $ws2.Cells.Item($rowIndex, 3).Value = $dep[2]
    $ws2.Cells.Item($rowIndex, 4).Value = $dep[3]
    $ws2.Cells.Item($rowIndex, 6).Value = $dep[4]
    $ws2.Cells.Item($rowIndex, 7).Value = $dep[5]
    $rowIndex++
}
